$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.205.46'
$ws.Range('D3').Value = '1.913.20'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8196'
$ws.Range('E5').Value = '  +4.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.14'
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3261'
$ws.Range('E8').Value = '  +3.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.06'
$ws.Range('E9').Value = '  +3.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07109'
$ws.Range('E10').Value = '  +2.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08083'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7805'
$ws.Range('E12').Value = '  +4.90%  '
$ws.Range('D13').Value = '1.926.41'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.343'
$ws.Range('E14').Value = '  +2.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.97'
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.35'
$ws.Range('E16').Value = '  +2.64%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '30.199.46'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.971'
$ws.Range('E18').Value = '  +1.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '249.20'
$ws.Range('E19').Value = '  +1.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007838'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.523'
$ws.Range('E23').Value = '  +9.36%  '
$ws.Range('B24').Value = 'Stellar'
$ws.Range('C24').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.1688'
$ws.Range('E24').Value = '  +22.78%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.455'
$ws.Range('E25').Value = '  +1.95%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.18'
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.08'
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.111'
$ws.Range('E28').Value = '  +3.89%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.371'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.533'
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.327'
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05731'
$ws.Range('E32').Value = '  +4.91%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.119'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.279'
$ws.Range('E34').Value = '  +1.77%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7397'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.726'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01936'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.806'
$ws.Range('E39').Value = '  +0.60%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4482'
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.76'
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.966'
$ws.Range('E42').Value = '  -2.78%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.928'
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8484'
$ws.Range('E44').Value = '  +1.51%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.93'
$ws.Range('E46').Value = '  +2.35%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.028.95'
$ws.Range('E47').Value = '  +4.93%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.912'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.604'
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.070.46'
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.561'
$ws.Range('E51').Value = '  +4.62%  '
